$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 325
$ws.Range("H2").Value = "kitchens"
$ws.Range("I2").Value = "distractor"
$ws.Range("K2").Value = "f"
$ws.Range("L2").Value = "stimuli/img_mucwi.png"
$ws.Range("M2").Value = 71.14814814814815
$ws.Range("N2").Value = 48.55555555555556
$ws.Range("O2").Value = 59.85185185185185
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 5

$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 326
$ws.Range("L3").Value = "stimuli/img_kost0.png"
$ws.Range("M3").Value = 63.09090909090909
$ws.Range("N3").Value = 42.77272727272727
$ws.Range("O3").Value = 52.93181818181819
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 5

$ws.Range("C4").Value = 1
$ws.Range("F4").Value = 327
$ws.Range("H4").Value = "bedrooms"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_badai.png"
$ws.Range("M4").Value = 63.97435897435897
$ws.Range("N4").Value = 43.38461538461539
$ws.Range("O4").Value = 53.67948717948718
$ws.Range("P4").Value = 39
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 4

$ws.Range("C5").Value = 1
$ws.Range("F5").Value = 328
$ws.Range("H5").Value = "kitchens"
$ws.Range("L5").Value = "stimuli/img_g7870.png"
$ws.Range("M5").Value = 68.70967741935483
$ws.Range("N5").Value = 44.2258064516129
$ws.Range("O5").Value = 56.46774193548387
$ws.Range("P5").Value = 31
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 4
$ws.Range("V5").Value = 4

$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 329
$ws.Range("H6").Value = "bedrooms"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_kn0we.png"
$ws.Range("M6").Value = 80.15909090909091
$ws.Range("N6").Value = 56.68181818181818
$ws.Range("O6").Value = 68.42045454545455
$ws.Range("P6").Value = 44
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 7
$ws.Range("T6").Value = 7
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = 7

$ws.Range("C7").Value = 1
$ws.Range("F7").Value = 330
$ws.Range("H7").Value = "living_rooms"
$ws.Range("I7").Value = "distractor"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_9z99v.png"
$ws.Range("M7").Value = 81.15625
$ws.Range("N7").Value = 64.78125
$ws.Range("O7").Value = 72.96875
$ws.Range("P7").Value = 32
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8
$ws.Range("S7").Value = 8
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 8
$ws.Range("V7").Value = 8

$ws.Range("C8").Value = 1
$ws.Range("F8").Value = 331
$ws.Range("H8").Value = "bedrooms"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_qbdgm.png"
$ws.Range("M8").Value = 76.88095238095238
$ws.Range("N8").Value = 60.40476190476191
$ws.Range("O8").Value = 68.64285714285714
$ws.Range("P8").Value = 42
$ws.Range("Q8").Value = 7
$ws.Range("R8").Value = 7
$ws.Range("S8").Value = 7
$ws.Range("T8").Value = 7
$ws.Range("U8").Value = 7
$ws.Range("V8").Value = 7

$ws.Range("C9").Value = 1
$ws.Range("F9").Value = 332
$ws.Range("H9").Value = "living_rooms"
$ws.Range("L9").Value = "stimuli/img_koooi.png"
$ws.Range("M9").Value = 63.95454545454545
$ws.Range("N9").Value = 44.56818181818182
$ws.Range("O9").Value = 54.26136363636364
$ws.Range("P9").Value = 44
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 5

$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 333
$ws.Range("H10").Value = "living_rooms"
$ws.Range("I10").Value = "distractor"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_6ddrx.png"
$ws.Range("M10").Value = 82.2
$ws.Range("N10").Value = 63.68571428571428
$ws.Range("O10").Value = 72.94285714285715
$ws.Range("P10").Value = 35
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 8
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = 8

$ws.Range("C11").Value = 1
$ws.Range("F11").Value = 334
$ws.Range("L11").Value = "stimuli/img_c2pbs.png"
$ws.Range("M11").Value = 21.95238095238095
$ws.Range("N11").Value = 14.47619047619048
$ws.Range("O11").Value = 18.21428571428572
$ws.Range("P11").Value = 42
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 1
$ws.Range("V11").Value = 1

$ws.Range("C12").Value = 1
$ws.Range("F12").Value = 335
$ws.Range("H12").Value = "living_rooms"
$ws.Range("I12").Value = "distractor"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_wz6x5.png"
$ws.Range("M12").Value = 68.3695652173913
$ws.Range("N12").Value = 48.47826086956522
$ws.Range("O12").Value = 58.42391304347826
$ws.Range("P12").Value = 46
$ws.Range("Q12").Value = 5
$ws.Range("R12").Value = 5
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 5
$ws.Range("U12").Value = 5
$ws.Range("V12").Value = 5

$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 336
$ws.Range("H13").Value = "living_rooms"
$ws.Range("I13").Value = "distractor"
$ws.Range("K13").Value = "f"
$ws.Range("L13").Value = "stimuli/img_o37la.png"
$ws.Range("M13").Value = 65.24324324324324
$ws.Range("N13").Value = 42.78378378378378
$ws.Range("O13").Value = 54.01351351351352
$ws.Range("P13").Value = 37
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 4
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 4
$ws.Range("U13").Value = 4
$ws.Range("V13").Value = 4

$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 337
$ws.Range("H14").Value = "bedrooms"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_jge7p.png"
$ws.Range("M14").Value = 90.42424242424242
$ws.Range("N14").Value = 75.63636363636364
$ws.Range("O14").Value = 83.03030303030303
$ws.Range("P14").Value = 33
$ws.Range("Q14").Value = 10
$ws.Range("R14").Value = 10
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 10
$ws.Range("U14").Value = 10
$ws.Range("V14").Value = 10

$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 338
$ws.Range("H15").Value = "living_rooms"
$ws.Range("I15").Value = "distractor"
$ws.Range("K15").Value = "f"
$ws.Range("L15").Value = "stimuli/img_xy930.png"
$ws.Range("M15").Value = 70.5952380952381
$ws.Range("N15").Value = 49.47619047619047
$ws.Range("O15").Value = 60.03571428571429
$ws.Range("P15").Value = 42
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("T15").Value = 5
$ws.Range("U15").Value = 5
$ws.Range("V15").Value = 5

$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 339
$ws.Range("H16").Value = "bedrooms"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_b971s.png"
$ws.Range("M16").Value = 70.5
$ws.Range("N16").Value = 47.61111111111111
$ws.Range("O16").Value = 59.05555555555556
$ws.Range("P16").Value = 36
$ws.Range("Q16").Value = 5
$ws.Range("R16").Value = 5
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 5
$ws.Range("U16").Value = 5
$ws.Range("V16").Value = 5

$ws.Range("C17").Value = 1
$ws.Range("F17").Value = 340
$ws.Range("H17").Value = "bedrooms"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_zv0dq.png"
$ws.Range("M17").Value = 76.86842105263158
$ws.Range("N17").Value = 52.71052631578947
$ws.Range("O17").Value = 64.78947368421052
$ws.Range("P17").Value = 38
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = 6
$ws.Range("S17").Value = 6
$ws.Range("T17").Value = 6
$ws.Range("U17").Value = 6
$ws.Range("V17").Value = 6

$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 341
$ws.Range("H18").Value = "kitchens"
$ws.Range("I18").Value = "distractor"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/img_d0k76.png"
$ws.Range("M18").Value = 67.09090909090909
$ws.Range("N18").Value = 46.3030303030303
$ws.Range("O18").Value = 56.6969696969697
$ws.Range("P18").Value = 33
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 4
$ws.Range("S18").Value = 4
$ws.Range("T18").Value = 4
$ws.Range("U18").Value = 4
$ws.Range("V18").Value = 4

$ws.Range("C19").Value = 1
$ws.Range("F19").Value = 342
$ws.Range("L19").Value = "stimuli/img_fbihy.png"
$ws.Range("M19").Value = 44.39024390243902
$ws.Range("N19").Value = 26.90243902439024
$ws.Range("O19").Value = 35.64634146341464
$ws.Range("P19").Value = 41
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 2
$ws.Range("U19").Value = 2
$ws.Range("V19").Value = 2

$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 343
$ws.Range("H20").Value = "living_rooms"
$ws.Range("L20").Value = "stimuli/img_53nbn.png"
$ws.Range("M20").Value = 73.28888888888889
$ws.Range("N20").Value = 51.15555555555556
$ws.Range("O20").Value = 62.22222222222223
$ws.Range("P20").Value = 45
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("T20").Value = 6
$ws.Range("U20").Value = 6
$ws.Range("V20").Value = 6

$ws.Range("C21").Value = 1
$ws.Range("F21").Value = 344
$ws.Range("H21").Value = "living_rooms"
$ws.Range("I21").Value = "distractor"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/img_hc49v.png"
$ws.Range("M21").Value = 70.95121951219512
$ws.Range("N21").Value = 53.31707317073171
$ws.Range("O21").Value = 62.13414634146342
$ws.Range("P21").Value = 41
$ws.Range("Q21").Value = 6
$ws.Range("R21").Value = 6
$ws.Range("S21").Value = 6
$ws.Range("T21").Value = 6
$ws.Range("U21").Value = 6
$ws.Range("V21").Value = 6

$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 345
$ws.Range("L22").Value = "stimuli/img_sfh4b.png"
$ws.Range("M22").Value = 69.06521739130434
$ws.Range("N22").Value = 49.54347826086956
$ws.Range("O22").Value = 59.30434782608695
$ws.Range("P22").Value = 46
$ws.Range("Q22").Value = 5
$ws.Range("R22").Value = 5
$ws.Range("S22").Value = 5
$ws.Range("T22").Value = 5
$ws.Range("U22").Value = 5
$ws.Range("V22").Value = 5

$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 346
$ws.Range("H23").Value = "living_rooms"
$ws.Range("L23").Value = "stimuli/img_73pyk.png"
$ws.Range("M23").Value = 69.27659574468085
$ws.Range("N23").Value = 47.27659574468085
$ws.Range("O23").Value = 58.27659574468085
$ws.Range("P23").Value = 47
$ws.Range("Q23").Value = 5
$ws.Range("R23").Value = 5
$ws.Range("S23").Value = 5
$ws.Range("T23").Value = 5
$ws.Range("U23").Value = 5
$ws.Range("V23").Value = 5

$ws.Range("C24").Value = 1
$ws.Range("F24").Value = 347
$ws.Range("H24").Value = "bedrooms"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_sltwe.png"
$ws.Range("M24").Value = 72.02500000000001
$ws.Range("N24").Value = 46.875
$ws.Range("O24").Value = 59.45
$ws.Range("P24").Value = 40
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 5

$ws.Range("C25").Value = 1
$ws.Range("F25").Value = 348
$ws.Range("H25").Value = "bedrooms"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_65cdi.png"
$ws.Range("M25").Value = 46.92307692307692
$ws.Range("N25").Value = 27
$ws.Range("O25").Value = 36.96153846153846
$ws.Range("P25").Value = 39
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 2
$ws.Range("T25").Value = 2
$ws.Range("U25").Value = 2
$ws.Range("V25").Value = 2

$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 349
$ws.Range("L26").Value = "stimuli/img_7os7q.png"
$ws.Range("M26").Value = 59.7027027027027
$ws.Range("N26").Value = 34.94594594594594
$ws.Range("O26").Value = 47.32432432432432
$ws.Range("P26").Value = 37

$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 350
$ws.Range("H27").Value = "bedrooms"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_i7vab.png"
$ws.Range("M27").Value = 86.40000000000001
$ws.Range("N27").Value = 67.8
$ws.Range("O27").Value = 77.09999999999999
$ws.Range("P27").Value = 35
$ws.Range("Q27").Value = 9
$ws.Range("R27").Value = 9
$ws.Range("S27").Value = 9
$ws.Range("T27").Value = 9
$ws.Range("U27").Value = 9
$ws.Range("V27").Value = 9

$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 351
$ws.Range("H28").Value = "living_rooms"
$ws.Range("I28").Value = "distractor"
$ws.Range("K28").Value = "f"
$ws.Range("L28").Value = "stimuli/img_swq34.png"
$ws.Range("M28").Value = 64.11363636363636
$ws.Range("N28").Value = 43.04545454545455
$ws.Range("O28").Value = 53.57954545454545
$ws.Range("P28").Value = 44
$ws.Range("Q28").Value = 5
$ws.Range("R28").Value = 5
$ws.Range("S28").Value = 5
$ws.Range("T28").Value = 5
$ws.Range("U28").Value = 5
$ws.Range("V28").Value = 5

$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 352
$ws.Range("H29").Value = "living_rooms"
$ws.Range("I29").Value = "distractor"
$ws.Range("K29").Value = "f"
$ws.Range("L29").Value = "stimuli/img_lgxzn.png"
$ws.Range("M29").Value = 73.11363636363636
$ws.Range("N29").Value = 49.97727272727273
$ws.Range("O29").Value = 61.54545454545455
$ws.Range("P29").Value = 44
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 5

$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 353
$ws.Range("H30").Value = "kitchens"
$ws.Range("L30").Value = "stimuli/img_pt3d7.png"
$ws.Range("M30").Value = 65.08571428571429
$ws.Range("N30").Value = 44.65714285714286
$ws.Range("O30").Value = 54.87142857142857
$ws.Range("P30").Value = 35
$ws.Range("Q30").Value = 4
$ws.Range("R30").Value = 4
$ws.Range("S30").Value = 4
$ws.Range("T30").Value = 4
$ws.Range("U30").Value = 4
$ws.Range("V30").Value = 4

$ws.Range("C31").Value = 1
$ws.Range("F31").Value = 354
$ws.Range("H31").Value = "bedrooms"
$ws.Range("I31").Value = "target"
$ws.Range("K31").Value = "j"
$ws.Range("L31").Value = "stimuli/img_bklr1.png"
$ws.Range("M31").Value = 86.54761904761905
$ws.Range("N31").Value = 67.73809523809524
$ws.Range("O31").Value = 77.14285714285714
$ws.Range("P31").Value = 42
$ws.Range("Q31").Value = 9
$ws.Range("R31").Value = 9
$ws.Range("S31").Value = 9
$ws.Range("T31").Value = 9
$ws.Range("U31").Value = 9
$ws.Range("V31").Value = 9

$ws.Range("C32").Value = 1
$ws.Range("F32").Value = 355
$ws.Range("H32").Value = "kitchens"
$ws.Range("I32").Value = "distractor"
$ws.Range("K32").Value = "f"
$ws.Range("L32").Value = "stimuli/img_4ufga.png"
$ws.Range("M32").Value = 67.79411764705883
$ws.Range("N32").Value = 41.5
$ws.Range("O32").Value = 54.64705882352941
$ws.Range("P32").Value = 34
$ws.Range("Q32").Value = 4
$ws.Range("R32").Value = 4
$ws.Range("S32").Value = 4
$ws.Range("T32").Value = 4
$ws.Range("U32").Value = 4
$ws.Range("V32").Value = 4

$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 356
$ws.Range("H33").Value = "bedrooms"
$ws.Range("I33").Value = "target"
$ws.Range("K33").Value = "j"
$ws.Range("L33").Value = "stimuli/img_bntrh.png"
$ws.Range("M33").Value = 76.07894736842105
$ws.Range("N33").Value = 53.36842105263158
$ws.Range("O33").Value = 64.72368421052632
$ws.Range("P33").Value = 38
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = 6
$ws.Range("S33").Value = 6
$ws.Range("T33").Value = 6
$ws.Range("U33").Value = 6
$ws.Range("V33").Value = 6

$ws.Range("C34").Value = 1
$ws.Range("F34").Value = 357
$ws.Range("H34").Value = "kitchens"
$ws.Range("I34").Value = "distractor"
$ws.Range("K34").Value = "f"
$ws.Range("L34").Value = "stimuli/img_oau79.png"
$ws.Range("M34").Value = 70.86486486486487
$ws.Range("N34").Value = 49
$ws.Range("O34").Value = 59.93243243243244
$ws.Range("P34").Value = 37

$ws.Range("C35").Value = 1
$ws.Range("F35").Value = 358
$ws.Range("H35").Value = "kitchens"
$ws.Range("I35").Value = "distractor"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_68wfw.png"
$ws.Range("M35").Value = 68.87878787878788
$ws.Range("N35").Value = 48.96969696969697
$ws.Range("O35").Value = 58.92424242424242
$ws.Range("P35").Value = 33
$ws.Range("Q35").Value = 5
$ws.Range("R35").Value = 5
$ws.Range("S35").Value = 5
$ws.Range("T35").Value = 5
$ws.Range("U35").Value = 5
$ws.Range("V35").Value = 5

$ws.Range("C36").Value = 1
$ws.Range("F36").Value = 359
$ws.Range("H36").Value = "bedrooms"
$ws.Range("I36").Value = "target"
$ws.Range("K36").Value = "j"
$ws.Range("L36").Value = "stimuli/img_dmjh8.png"
$ws.Range("M36").Value = 57.48648648648648
$ws.Range("N36").Value = 37.64864864864865
$ws.Range("O36").Value = 47.56756756756756
$ws.Range("P36").Value = 37
$ws.Range("Q36").Value = 3
$ws.Range("R36").Value = 3
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 3
$ws.Range("U36").Value = 3
$ws.Range("V36").Value = 3

$ws.Range("C37").Value = 1
$ws.Range("F37").Value = 360
$ws.Range("L37").Value = "stimuli/img_ri0yx.png"
$ws.Range("M37").Value = 88.96969696969697
$ws.Range("N37").Value = 77.15151515151516
$ws.Range("O37").Value = 83.06060606060606
$ws.Range("P37").Value = 33
$ws.Range("Q37").Value = 10
$ws.Range("R37").Value = 10
$ws.Range("S37").Value = 10
$ws.Range("T37").Value = 10
$ws.Range("U37").Value = 10
$ws.Range("V37").Value = 10

$ws.Range("C38").Value = 1
$ws.Range("F38").Value = 361
$ws.Range("H38").Value = "kitchens"
$ws.Range("I38").Value = "distractor"
$ws.Range("K38").Value = "f"
$ws.Range("L38").Value = "stimuli/img_anjr0.png"
$ws.Range("M38").Value = 67.88888888888889
$ws.Range("N38").Value = 45.80555555555556
$ws.Range("O38").Value = 56.84722222222222
$ws.Range("P38").Value = 36
$ws.Range("Q38").Value = 4
$ws.Range("R38").Value = 4
$ws.Range("S38").Value = 4
$ws.Range("T38").Value = 4
$ws.Range("U38").Value = 4
$ws.Range("V38").Value = 4

$ws.Range("C39").Value = 1
$ws.Range("F39").Value = 362
$ws.Range("L39").Value = "stimuli/img_l1h36.png"
$ws.Range("M39").Value = 26.64285714285714
$ws.Range("N39").Value = 9.142857142857142
$ws.Range("O39").Value = 17.89285714285714
$ws.Range("P39").Value = 42
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = 1
$ws.Range("S39").Value = 1
$ws.Range("T39").Value = 1
$ws.Range("U39").Value = 1
$ws.Range("V39").Value = 1

$ws.Range("C40").Value = 1
$ws.Range("F40").Value = 363
$ws.Range("H40").Value = "kitchens"
$ws.Range("L40").Value = "stimuli/img_q1ynd.png"
$ws.Range("M40").Value = 70.05714285714286
$ws.Range("N40").Value = 47.31428571428572
$ws.Range("O40").Value = 58.68571428571429
$ws.Range("P40").Value = 35
$ws.Range("Q40").Value = 5
$ws.Range("R40").Value = 5
$ws.Range("S40").Value = 5
$ws.Range("T40").Value = 5
$ws.Range("U40").Value = 5
$ws.Range("V40").Value = 5

$ws.Range("C41").Value = 1
$ws.Range("F41").Value = 364
$ws.Range("H41").Value = "bedrooms"
$ws.Range("I41").Value = "target"
$ws.Range("K41").Value = "j"
$ws.Range("L41").Value = "stimuli/img_5mw7y.png"
$ws.Range("M41").Value = 72.65909090909091
$ws.Range("N41").Value = 50.86363636363637
$ws.Range("O41").Value = 61.76136363636364
$ws.Range("P41").Value = 44
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = 6
$ws.Range("S41").Value = 6
$ws.Range("T41").Value = 6
$ws.Range("U41").Value = 6
$ws.Range("V41").Value = 6

